# Generate Report for handoff
# Adds two new tracked files (03acf302-...md and ac6c1e5d-...md) to the
# localization status report, pushes the existing ".localization-config"
# row down, and flips the two previously-tracked files' status from
# "Ready for handoff" to "In Translation" on the language sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$mdCommit = "c174fd49aced88769ee5a2cb675486f334333fc9"
$zhCommit = "d25626040259023ca675ca378693f1f437e3b53a"
$deCommit = "111d40962d8a64ea9b1426333d6970210bd7e572"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/.localization-config"
$zhBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/"
$deBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/"

$file3 = "03acf302-3652-4579-973c-b89bbd18ffc4.md"
$file3ZhXlf = "03acf302-3652-4579-973c-b89bbd18ffc4.a1b011afc322fa2ea2e491f688a206adf34c2fc5.zh-cn.xlf"
$file3DeXlf = "03acf302-3652-4579-973c-b89bbd18ffc4.a1b011afc322fa2ea2e491f688a206adf34c2fc5.de-de.xlf"

$file4 = "ac6c1e5d-c464-49b4-8cda-02dce0384b70.md"
$file4ZhXlf = "ac6c1e5d-c464-49b4-8cda-02dce0384b70.b44c68e17ae6ad15d3d3f0f0e344be729f4e1235.zh-cn.xlf"
$file4DeXlf = "ac6c1e5d-c464-49b4-8cda-02dce0384b70.b44c68e17ae6ad15d3d3f0f0e344be729f4e1235.de-de.xlf"

# ---------------------------------------------------------------------------
# 1. Overview sheet: shift the ".localization-config" summary row from row 4
#    to row 6, and insert the two new files as "Ready for handoff" rows 4-5.
# ---------------------------------------------------------------------------

$wsOverview.Range("A6").Value = ".localization-config"
$wsOverview.Range("B6").Value = "Not to be localized"
$wsOverview.Range("C6").Value = "Not to be localized"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $mdBase + $file3, "", "", $file3)
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), $mdBase + $file4, "", "", $file4)
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), $cfgUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: rows 2-3 move to "In Translation"; ".localization-config"
#    shifts from row 4 to row 6; the two new files become rows 4-5.
# ---------------------------------------------------------------------------

$wsZh.Range("B2").Value = "In Translation"
$wsZh.Range("B3").Value = "In Translation"

$wsZh.Range("A6").Value = ".localization-config"
$wsZh.Range("B6").Value = "Not to be localized"
$wsZh.Range("D6").Value = "0001-01-01 00:00:00"
$wsZh.Range("G6").Value = "0001-01-01 00:00:00"
$wsZh.Range("H6").Value = "Ignored"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $mdBase + $file3, "", "", $file3)
$wsZh.Range("B4").Value = "Ready for handoff"
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), $zhBase + $file3ZhXlf, "", "", $file3ZhXlf)
$wsZh.Range("D4").Value = "2016-01-25 05:35:29"
$wsZh.Range("G4").Value = "0001-01-01 00:00:00"
$wsZh.Range("H4").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), $mdBase + $file4, "", "", $file4)
$wsZh.Range("B5").Value = "Ready for handoff"
$wsZh.Hyperlinks.Add($wsZh.Range("C5"), $zhBase + $file4ZhXlf, "", "", $file4ZhXlf)
$wsZh.Range("D5").Value = "2016-01-25 05:35:29"
$wsZh.Range("G5").Value = "0001-01-01 00:00:00"
$wsZh.Range("H5").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A6"), $cfgUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# 3. de-de sheet: same shape as zh-cn, with the de-de specific file names.
# ---------------------------------------------------------------------------

$wsDe.Range("B2").Value = "In Translation"
$wsDe.Range("B3").Value = "In Translation"

$wsDe.Range("A6").Value = ".localization-config"
$wsDe.Range("B6").Value = "Not to be localized"
$wsDe.Range("D6").Value = "0001-01-01 00:00:00"
$wsDe.Range("G6").Value = "0001-01-01 00:00:00"
$wsDe.Range("H6").Value = "Ignored"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $mdBase + $file3, "", "", $file3)
$wsDe.Range("B4").Value = "Ready for handoff"
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), $deBase + $file3DeXlf, "", "", $file3DeXlf)
$wsDe.Range("D4").Value = "2016-01-25 05:35:39"
$wsDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDe.Range("H4").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $mdBase + $file4, "", "", $file4)
$wsDe.Range("B5").Value = "Ready for handoff"
$wsDe.Hyperlinks.Add($wsDe.Range("C5"), $deBase + $file4DeXlf, "", "", $file4DeXlf)
$wsDe.Range("D5").Value = "2016-01-25 05:35:39"
$wsDe.Range("G5").Value = "0001-01-01 00:00:00"
$wsDe.Range("H5").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A6"), $cfgUrl, "", "", ".localization-config")

Write-Host "Localization status report updated."
